$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4 values (mirrors structure of row 3)
$ws.Range("A4").Value = 42641.540729166663
$ws.Range("B4").Value = $true
$ws.Range("C4").Value = 9976.89
$ws.Range("D4").Value = 9953
$ws.Range("E4").Value = 79.319999999999993
$ws.Range("F4").Value = 78.94
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = -0.48
$ws.Range("I4").Value = $false

# Copy date-style formatting from row 3's date columns (A and G) onto row 4
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)
